$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 15.0
$ws.Range("B16").Value = "dealAnalysis_Finance_231310_TC_02"
$ws.Range("D16").Value = "25/04/2022"
$ws.Range("E16").Value = "Fail"
